# [TEST][OWNBRANCH] Slack and git Notifications Testing on my own branch v2
#
# Append, at the end of the document, an empty paragraph followed by a
# paragraph containing "Test 2 ". Both new paragraphs carry the same
# German (de-DE) language mark already used by the existing text, so the
# run/paragraph properties are an exact match for the rest of the file.

$d = $word.ActiveDocument

# Collapse a range to the very end of the story so the new content lands
# right before the final section properties, just like typing at the end
# of the document in the Word UI.
$endOfDoc = $d.Content
$endOfDoc.Collapse(0)

# Build the two new paragraphs as a WordprocessingML package fragment so
# the empty paragraph really stays empty (no stray run is created) and the
# "Test 2 " paragraph gets its own run with identical formatting.
$newParagraphsXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:lang w:val="de-DE"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:lang w:val="de-DE"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:val="de-DE"/>
              </w:rPr>
              <w:t xml:space="preserve">Test 2 </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$endOfDoc.InsertXML($newParagraphsXml) | Out-Null
